# Sample Project / Main.xlsx — "Rules" sheet
#
# The rule row that was labelled "R40" (cell B11 of the decision table)
# is renamed to "1".
#
# A leading apostrophe is used so Excel stores the new value as TEXT
# (matching the original cell's string type, t="s") instead of silently
# auto-converting the digit "1" into a numeric value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = "'1"
